# Auto-generated script applying cryptos.xlsx price/volume update diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to keep its original text type even if the
    # string looks like a number (e.g. "0.0620" or "214.47"),
    # by using an apostrophe-prefixed literal, then restore the
    # original cell style so no stray quote-prefix formatting remains.
    $origStyle = $range.Style
    $range.Value = "'" + $text
    $range.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '27.587.20'
$ws.Range('E2').Value = '  +1.54%  '
Set-TextValue $ws.Range('D3') '1.656.45'
$ws.Range('E3').Value = '  -1.27%  '
$ws.Range('E4').Value = '  +0.33%  '
Set-TextValue $ws.Range('D5') '214.47'
$ws.Range('E5').Value = '  +0.03%  '
$ws.Range('E6').Value = '  -0.91%  '
$ws.Range('E7').Value = '  +0.40%  '
Set-TextValue $ws.Range('D8') '23.27'
$ws.Range('E8').Value = '  +1.33%  '
$ws.Range('E9').Value = '  -0.77%  '
Set-TextValue $ws.Range('D10') '0.0620'
$ws.Range('E10').Value = '  -0.37%  '
Set-TextValue $ws.Range('D11') '0.0879'
$ws.Range('E11').Value = '  -1.32%  '
Set-TextValue $ws.Range('D12') '1.891.83'
$ws.Range('E12').Value = '  -1.23%  '
Set-TextValue $ws.Range('D13') '1.658.80'
$ws.Range('E13').Value = '  -1.15%  '
$ws.Range('E14').Value = '  -1.91%  '
Set-TextValue $ws.Range('D15') '0.548'
$ws.Range('E15').Value = '  -2.02%  '
Set-TextValue $ws.Range('D16') '65.76'
$ws.Range('E16').Value = '  -1.34%  '
Set-TextValue $ws.Range('D17') '246.84'
$ws.Range('E17').Value = '  +4.81%  '
Set-TextValue $ws.Range('D18') '27.613.34'
$ws.Range('E18').Value = '  +1.73%  '
Set-TextValue $ws.Range('D19') '0.0₃0728'
$ws.Range('E19').Value = '  -1.91%  '
Set-TextValue $ws.Range('D20') '7.44'
$ws.Range('E20').Value = '  -5.03%  '
$ws.Range('E21').Value = '  +0.22%  '
Set-TextValue $ws.Range('D22') '4.46'
$ws.Range('E22').Value = '  -2.09%  '
Set-TextValue $ws.Range('D23') '9.07'
$ws.Range('E23').Value = '  -5.01%  '
Set-TextValue $ws.Range('D24') '2.02'
$ws.Range('E24').Value = '  -3.52%  '
Set-TextValue $ws.Range('D25') '146.42'
$ws.Range('E25').Value = '  -0.88%  '
Set-TextValue $ws.Range('D26') '7.16'
$ws.Range('E26').Value = '  -4.46%  '
$ws.Range('E28').Value = '  +0.35%  '
$ws.Range('E29').Value = '  -1.41%  '
$ws.Range('E30').Value = '  +6.06%  '
$ws.Range('E31').Value = '  -0.15%  '
Set-TextValue $ws.Range('D32') '3.34'
$ws.Range('E32').Value = '  -1.08%  '
Set-TextValue $ws.Range('D33') '1.429.03'
$ws.Range('E33').Value = '  -7.05%  '
Set-TextValue $ws.Range('D34') '3.12'
$ws.Range('E34').Value = '  -3.82%  '
$ws.Range('E35').Value = '  -7.66%  '
Set-TextValue $ws.Range('D36') '2.39'
$ws.Range('E36').Value = '  +0.38%  '
Set-TextValue $ws.Range('D37') '0.926'
$ws.Range('E37').Value = '  -2.54%  '
$ws.Range('E38').Value = '  -5.02%  '
Set-TextValue $ws.Range('D40') '1.04'
$ws.Range('E40').Value = '  -2.10%  '
Set-TextValue $ws.Range('D41') '69.15'
$ws.Range('E41').Value = '  -1.06%  '
$ws.Range('E42').Value = '  +0.29%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D43') '5.39'
$ws.Range('E43').Value = '  -6.81%  '
$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range('D44') '2.21'
$ws.Range('E44').Value = '  -1.73%  '
Set-TextValue $ws.Range('D45') '1.799.73'
$ws.Range('E45').Value = '  -1.41%  '
Set-TextValue $ws.Range('D46') '0.787'
$ws.Range('E46').Value = '  +0.64%  '
$ws.Range('E47').Value = '  +2.73%  '
Set-TextValue $ws.Range('D48') '88.63'
$ws.Range('E48').Value = '  -1.71%  '
$ws.Range('E49').Value = '  -4.27%  '
$ws.Range('E50').Value = '  -3.66%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D51') '7.77'
$ws.Range('E51').Value = '  -5.36%  '
